# Asset Page detail changes + test data updates to select the lookup of
# contact name (per commit message). Updates several "QA_Automation" test
# data sheets: Campaign, Account, Contact, Asset, Opportunity, Quote, and
# repositions the active selection on several sheets. Contract becomes the
# final active/selected tab.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# LoginPage: selection moves from D6 to C2 (no longer the "active" tab)
# ---------------------------------------------------------------------
$wsLogin = $wb.Worksheets.Item("LoginPage")
$wsLogin.Activate()
$wsLogin.Range("C2").Select()

# ---------------------------------------------------------------------
# Campaign: campaign name test data updated
# ---------------------------------------------------------------------
$wsCampaign = $wb.Worksheets.Item("Campaign")
$wsCampaign.Activate()
$wsCampaign.Range("A2").Value = "Test Campaign QA"
$wsCampaign.Range("A8").Select()

# ---------------------------------------------------------------------
# Account: excel input test data updated, selection moves to B4,
# the window no longer needs to be scrolled to column D
# ---------------------------------------------------------------------
$wsAccount = $wb.Worksheets.Item("Account")
$wsAccount.Activate()
$wsAccount.Range("A2").Value = "Test Account excel input"
$wsAccount.Range("B4").Select()

# ---------------------------------------------------------------------
# Contact: updated to pick a lookup contact ("Test QA" / "Contact" /
# "Test Asset Account"), selection moves to D9
# ---------------------------------------------------------------------
$wsContact = $wb.Worksheets.Item("Contact")
$wsContact.Activate()
$wsContact.Range("A2").Value = "Test QA"
$wsContact.Range("B2").Value = "Contact"
$wsContact.Range("C2").Value = "Test Asset Account"
$wsContact.Columns.Item(1).ColumnWidth = 8.893229166666666
$wsContact.Columns.Item(2).ColumnWidth = 8.619791666666666
$wsContact.Columns.Item(3).ColumnWidth = 15.619791666666666
$wsContact.Range("D9").Select()

# ---------------------------------------------------------------------
# Asset: new "UpdatedAssetName" column (D) added with highlighter data,
# existing row renamed to the new QA asset names, selection moves to B2
# ---------------------------------------------------------------------
$wsAsset = $wb.Worksheets.Item("Asset")
$wsAsset.Activate()
$wsAsset.Range("A2").Value = "Test QA Asset"
$wsAsset.Range("B2").Value = "Test Asset Account"
$wsAsset.Range("C2").Value = "Test Asset Contact"
$wsAsset.Range("D1").Value = "UpdatedAssetName"
$wsAsset.Range("D2").Value = "Test QA Asset Updated"
$wsAsset.Columns.Item(1).ColumnWidth = 11.346354166666666
$wsAsset.Columns.Item(2).ColumnWidth = 15.619791666666666
$wsAsset.Columns.Item(4).ColumnWidth = 19.166666666666668
$wsAsset.Range("B2").Select()

# ---------------------------------------------------------------------
# Opportunity: sample opportunity name updated, selection moves to A3
# ---------------------------------------------------------------------
$wsOpportunity = $wb.Worksheets.Item("Opportunity")
$wsOpportunity.Activate()
$wsOpportunity.Range("A2").Value = "Test My Opportunity"
$wsOpportunity.Range("A3").Select()

# ---------------------------------------------------------------------
# Quote: sample quote name updated, selection moves to A3
# ---------------------------------------------------------------------
$wsQuote = $wb.Worksheets.Item("Quote")
$wsQuote.Activate()
$wsQuote.Range("A2").Value = "Test My Quote"
$wsQuote.Range("A3").Select()

# ---------------------------------------------------------------------
# Contract: becomes the active / selected tab, selection moves to C2
# ---------------------------------------------------------------------
$wsContract = $wb.Worksheets.Item("Contract")
$wsContract.Activate()
$wsContract.Range("C2").Select()
